$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.004.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "'2.358.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'0.680"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").Value = "'239.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("D7").Value = "'74.23"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.34%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.597"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +11.11%  "
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").Value = "'57.27"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "'32.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.19%  "
$ws.Range("E13").Value = "  +9.63%  "
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").Value = "'2.710.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "'16.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("D17").Value = "'0.899"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").Value = "'2.361.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "'43.903.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").Value = "'6.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.95%  "
$ws.Range("D22").Value = "'76.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("D23").Value = "'258.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.28%  "
$ws.Range("D24").Value = "'1.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +24.56%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("D27").Value = "'3.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.26%  "
$ws.Range("D28").Value = "'10.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.11%  "
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("D30").Value = "'22.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.44%  "
$ws.Range("D31").Value = "'175.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.68%  "
$ws.Range("D32").Value = "'0.128"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.17%  "
$ws.Range("D33").Value = "'0.136"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.21%  "
$ws.Range("D34").Value = "'0.0772"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.35%  "
$ws.Range("D35").Value = "'5.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.02%  "
$ws.Range("D36").Value = "'5.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.69%  "
$ws.Range("D37").Value = "'3.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.31%  "
$ws.Range("E38").Value = "  -3.17%  "
$ws.Range("D39").Value = "'6.31"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("E40").Value = "  +4.81%  "
$ws.Range("D41").Value = "'0.112"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.57%  "
$ws.Range("D42").Value = "'0.207"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +14.55%  "
$ws.Range("D43").Value = "'9.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.04%  "
$ws.Range("D44").Value = "'19.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.01%  "
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").Value = "'4.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.38%  "
$ws.Range("D47").Value = "'2.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.95%  "
$ws.Range("D48").Value = "'58.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.93%  "
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("D51").Value = "'100.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.89%  "
